$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.744.71'

$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '1.646.43'

$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.19'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  +1.15%  '

$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = '  +0.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = '  +2.04%  '

$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").Value = '1.877.88'

$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("D13").Value = '1.642.45'

$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  +1.16%  '

$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.21'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("D17").Value = '26.754.06'

$ws.Range("E17").Value = '  +0.67%  '

$ws.Range("E18").Value = '  +0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.16'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.54'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = '  +14.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = '  +1.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.27'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("E24").Value = '  +1.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.49'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.77'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  +0.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = '  +1.26%  '

$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = '  -0.32%  '

$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("D34").Value = '1.279.44'

$ws.Range("E34").Value = '  +1.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = '  +2.40%  '

$ws.Range("E36").Value = '  +2.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0179'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  +2.02%  '

$ws.Range("E38").Value = '  +5.81%  '

$ws.Range("E39").Value = '  +3.98%  '

$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.815'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.24'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  -1.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.44'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  +1.60%  '

$ws.Range("D44").Value = '1.789.33'

$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.94'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  -1.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.90'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  +8.76%  '

$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("B48").Value = 'Cronos'

$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0515'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("B49").Value = 'BabyDogeCoin'

$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'

$ws.Range("D49").Value = '0.0₆0100'

$ws.Range("E49").Value = '  -3.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.77'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  +2.13%  '
